$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6947294.5
$ws.Range("I76").Value = 3160
$ws.Range("J76").Value = 18520852
$ws.Range("K76").Value = 3160
$ws.Range("L76").Value = 18520852
$ws.Range("M76").Value = -2845
$ws.Range("N76").Value = -18521482

$ws.Range("H79").Value = 6947294.5
$ws.Range("I79").Value = 3160
$ws.Range("J79").Value = 18520852
$ws.Range("K79").Value = 3160
$ws.Range("L79").Value = 18520852
$ws.Range("M79").Value = -2068
$ws.Range("N79").Value = -18523036

$ws.Range("H129").Value = 859.725
$ws.Range("J129").Value = 874.7105
$ws.Range("L129").Value = 2624.1315
$ws.Range("N129").Value = -12624.1315

$ws.Range("H132").Value = 34942.97
$ws.Range("I132").Value = 44726.75
$ws.Range("K132").Value = 134180.25
$ws.Range("M132").Value = -131650.25

$ws.Range("H133").Value = 51777.5
$ws.Range("J133").Value = 51777.5
$ws.Range("L133").Value = 51777.5
$ws.Range("N133").Value = -61897.5

$ws.Range("H138").Value = 2877.25
$ws.Range("J138").Value = 3096.9375
$ws.Range("L138").Value = 9290.8125
$ws.Range("N138").Value = -19570.8125


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3123.923
$ws.Range("I45").Value = 2245.3333
$ws.Range("K45").Value = 2245.3333
$ws.Range("M45").Value = -1868.3333

$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws.Range("H114").Value = 43066.332
$ws.Range("J114").Value = 43066.332
$ws.Range("L114").Value = 43066.332
$ws.Range("N114").Value = -51744.332

$ws.Range("H124").Value = 11596.625
$ws.Range("J124").Value = 11596.625
$ws.Range("L124").Value = 11596.625
$ws.Range("N124").Value = -21416.625

$ws.Range("H125").Value = 36995
$ws.Range("J125").Value = 36995
$ws.Range("L125").Value = 36995
$ws.Range("N125").Value = -46835

$ws.Range("H132").Value = 17993.838
$ws.Range("I132").Value = 1554.6111
$ws.Range("K132").Value = 4663.8333
$ws.Range("M132").Value = -2133.8333


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13697.027
$ws.Range("I31").Value = 21197.55
$ws.Range("K31").Value = 21197.55
$ws.Range("M31").Value = -20902.55

$ws.Range("H34").Value = 13697.027
$ws.Range("I34").Value = 21197.55
$ws.Range("K34").Value = 21197.55
$ws.Range("M34").Value = -20995.55

$ws.Range("H122").Value = 2519.4546
$ws.Range("I122").Value = 2726.75
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 8180.25
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -5730.25
$ws.Range("N122").Value = -10799.9998


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8983.308000000001
$ws.Range("J68").Value = 10494
$ws.Range("L68").Value = 31482
$ws.Range("N68").Value = -33104

$ws.Range("H71").Value = 8983.308000000001
$ws.Range("J71").Value = 10494
$ws.Range("L71").Value = 94446
$ws.Range("N71").Value = -102558

$ws.Range("H131").Value = 114454.59
$ws.Range("I131").Value = 757.5
$ws.Range("J131").Value = 125824.3
$ws.Range("K131").Value = 2272.5
$ws.Range("L131").Value = 377472.9
$ws.Range("M131").Value = 2767.5
$ws.Range("N131").Value = -387552.9

$ws.Range("H136").Value = 1999.1875
$ws.Range("I136").Value = 999.1667
$ws.Range("J136").Value = 4999.25
$ws.Range("K136").Value = 2997.5001
$ws.Range("L136").Value = 14997.75
$ws.Range("M136").Value = 2102.4999
$ws.Range("N136").Value = -25197.75

$ws.Range("H137").Value = 18522360
$ws.Range("J137").Value = 23814118
$ws.Range("L137").Value = 71442354
$ws.Range("N137").Value = -71452554

$ws.Range("H138").Value = 1923.1177
$ws.Range("I138").Value = 1622.7273
$ws.Range("J138").Value = 2473.8333
$ws.Range("K138").Value = 4868.1819
$ws.Range("L138").Value = 7421.499899999999
$ws.Range("M138").Value = 271.8181000000004
$ws.Range("N138").Value = -17701.4999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9282.588
$ws.Range("I80").Value = 17022
$ws.Range("J80").Value = 3865
$ws.Range("K80").Value = 17022
$ws.Range("L80").Value = 3865
$ws.Range("M80").Value = -16024
$ws.Range("N80").Value = -5861

$ws.Range("H83").Value = 9282.588
$ws.Range("I83").Value = 17022
$ws.Range("J83").Value = 3865
$ws.Range("K83").Value = 85110
$ws.Range("L83").Value = 19325
$ws.Range("M83").Value = -80118
$ws.Range("N83").Value = -29309

$ws.Range("H132").Value = 21228.928
$ws.Range("J132").Value = 48380.637
$ws.Range("L132").Value = 145141.911
$ws.Range("N132").Value = -150201.911

$ws.Range("H134").Value = 29576
$ws.Range("J134").Value = 29576
$ws.Range("L134").Value = 88728
$ws.Range("N134").Value = -93798


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3662.5
$ws.Range("I100").Value = 1550
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 1550
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -1009
$ws.Range("N100").Value = -11082

$ws.Range("H110").Value = 2529725
$ws.Range("J110").Value = 2529725
$ws.Range("L110").Value = 2529725
$ws.Range("N110").Value = -2537905

$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920

$ws.Range("H132").Value = 1543.7297
$ws.Range("I132").Value = 1051.9354
$ws.Range("J132").Value = 4084.6667
$ws.Range("K132").Value = 3155.8062
$ws.Range("L132").Value = 12254.0001
$ws.Range("M132").Value = -625.8062
$ws.Range("N132").Value = -17314.0001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5400
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5400
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5400
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6648

$ws.Range("H65").Value = 5400
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5400
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 27000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -33240

$ws.Range("H113").Value = 2079987.4
$ws.Range("I113").Value = 1328.3334
$ws.Range("K113").Value = 3985.0002
$ws.Range("M113").Value = -1815.0002

$ws.Range("H132").Value = 2124.814
$ws.Range("I132").Value = 1874
$ws.Range("K132").Value = 5622
$ws.Range("M132").Value = -3092

